$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17. This shifts the existing rows 17-23
# down to 18-24, preserving all of their data (matching the diff, where
# the old row 17-23 data reappears unchanged one row lower).
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly record.
$ws.Cells.Item(17, 1).Value = 7
$ws.Cells.Item(17, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(17, 3).Value = "Ñuble"
$ws.Cells.Item(17, 4).Value = 44609
$ws.Cells.Item(17, 5).Value = 16
$ws.Cells.Item(17, 6).Value = 100112040
$ws.Cells.Item(17, 7).Value = "Cilantro"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 60
$ws.Cells.Item(17, 11).Value = 550
$ws.Cells.Item(17, 12).Value = 600
$ws.Cells.Item(17, 13).Value = 575
$ws.Cells.Item(17, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(17, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(17, 16).Value = 575
$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = "Hortaliza"

# Make sure the new date cell keeps the same date number format as the
# rest of column D (the Insert() above should already have copied it,
# but set it explicitly to be safe).
$ws.Cells.Item(17, 4).NumberFormat = $ws.Cells.Item(18, 4).NumberFormat
